# Auto-generated Excel COM-interop script: applies the scheduled market-data refresh
# to the leve-profit tables (columns H-N) on each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1050.25
$ws.Cells.Item(40, 9).Value = 1000
$ws.Cells.Item(40, 10).Value = 1057.4286
$ws.Cells.Item(40, 11).Value = 1000
$ws.Cells.Item(40, 12).Value = 1057.4286
$ws.Cells.Item(40, 13).Value = -825
$ws.Cells.Item(40, 14).Value = -1407.4286
$ws.Cells.Item(64, 8).Value = 4222.864
$ws.Cells.Item(64, 10).Value = 5422.143
$ws.Cells.Item(64, 12).Value = 5422.143
$ws.Cells.Item(64, 14).Value = -5918.143
$ws.Cells.Item(67, 8).Value = 4222.864
$ws.Cells.Item(67, 10).Value = 5422.143
$ws.Cells.Item(67, 12).Value = 5422.143
$ws.Cells.Item(67, 14).Value = -7138.143
$ws.Cells.Item(76, 8).Value = 7343.5586
$ws.Cells.Item(76, 9).Value = 10073.8125
$ws.Cells.Item(76, 10).Value = 4916.6665
$ws.Cells.Item(76, 11).Value = 10073.8125
$ws.Cells.Item(76, 12).Value = 4916.6665
$ws.Cells.Item(76, 13).Value = -9758.8125
$ws.Cells.Item(76, 14).Value = -5546.6665
$ws.Cells.Item(79, 8).Value = 7343.5586
$ws.Cells.Item(79, 9).Value = 10073.8125
$ws.Cells.Item(79, 10).Value = 4916.6665
$ws.Cells.Item(79, 11).Value = 10073.8125
$ws.Cells.Item(79, 12).Value = 4916.6665
$ws.Cells.Item(79, 13).Value = -8981.8125
$ws.Cells.Item(79, 14).Value = -7100.6665
$ws.Cells.Item(88, 8).Value = 3127.6956
$ws.Cells.Item(88, 9).Value = 480.57144
$ws.Cells.Item(88, 10).Value = 4285.8125
$ws.Cells.Item(88, 11).Value = 480.57144
$ws.Cells.Item(88, 12).Value = 4285.8125
$ws.Cells.Item(88, 13).Value = -74.57144
$ws.Cells.Item(88, 14).Value = -5097.8125
$ws.Cells.Item(91, 8).Value = 3127.6956
$ws.Cells.Item(91, 9).Value = 480.57144
$ws.Cells.Item(91, 10).Value = 4285.8125
$ws.Cells.Item(91, 11).Value = 480.57144
$ws.Cells.Item(91, 12).Value = 4285.8125
$ws.Cells.Item(91, 13).Value = 923.4285600000001
$ws.Cells.Item(91, 14).Value = -7093.8125
$ws.Cells.Item(112, 8).Value = 1911.2963
$ws.Cells.Item(112, 9).Value = 995
$ws.Cells.Item(112, 10).Value = 1946.5385
$ws.Cells.Item(112, 11).Value = 2985
$ws.Cells.Item(112, 12).Value = 5839.6155
$ws.Cells.Item(112, 13).Value = -1877
$ws.Cells.Item(112, 14).Value = -8055.6155
$ws.Cells.Item(116, 8).Value = 127440.3
$ws.Cells.Item(116, 9).Value = 164691.16
$ws.Cells.Item(116, 10).Value = 6375
$ws.Cells.Item(116, 11).Value = 164691.16
$ws.Cells.Item(116, 12).Value = 6375
$ws.Cells.Item(116, 13).Value = -161249.16
$ws.Cells.Item(116, 14).Value = -13259
$ws.Cells.Item(126, 8).Value = 39950
$ws.Cells.Item(126, 10).Value = 39950
$ws.Cells.Item(126, 12).Value = 39950
$ws.Cells.Item(126, 14).Value = -49830
$ws.Cells.Item(132, 8).Value = 3468.82
$ws.Cells.Item(132, 9).Value = 1698.6888
$ws.Cells.Item(132, 10).Value = 19400
$ws.Cells.Item(132, 11).Value = 5096.0664
$ws.Cells.Item(132, 12).Value = 58200
$ws.Cells.Item(132, 13).Value = -2566.0664
$ws.Cells.Item(132, 14).Value = -63260
$ws.Cells.Item(135, 8).Value = 443.2258
$ws.Cells.Item(135, 9).Value = 383.60715
$ws.Cells.Item(135, 10).Value = 999.6667
$ws.Cells.Item(135, 11).Value = 3452.46435
$ws.Cells.Item(135, 12).Value = 8997.0003
$ws.Cells.Item(135, 13).Value = -917.4643499999997
$ws.Cells.Item(135, 14).Value = -14067.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 728.5833
$ws.Cells.Item(2, 9).Value = 660.3333
$ws.Cells.Item(2, 11).Value = 660.3333
$ws.Cells.Item(2, 13).Value = -547.3333
$ws.Cells.Item(63, 8).Value = 6708.7896
$ws.Cells.Item(63, 9).Value = 7713.9165
$ws.Cells.Item(63, 10).Value = 4985.7144
$ws.Cells.Item(63, 11).Value = 7713.9165
$ws.Cells.Item(63, 12).Value = 4985.7144
$ws.Cells.Item(63, 13).Value = -7027.9165
$ws.Cells.Item(63, 14).Value = -6357.7144
$ws.Cells.Item(66, 8).Value = 6708.7896
$ws.Cells.Item(66, 9).Value = 7713.9165
$ws.Cells.Item(66, 10).Value = 4985.7144
$ws.Cells.Item(66, 11).Value = 38569.5825
$ws.Cells.Item(66, 12).Value = 24928.572
$ws.Cells.Item(66, 13).Value = -35137.5825
$ws.Cells.Item(66, 14).Value = -31792.572
$ws.Cells.Item(87, 8).Value = 33000
$ws.Cells.Item(87, 10).Value = 33000
$ws.Cells.Item(87, 12).Value = 33000
$ws.Cells.Item(87, 14).Value = -35496
$ws.Cells.Item(90, 8).Value = 33000
$ws.Cells.Item(90, 10).Value = 33000
$ws.Cells.Item(90, 12).Value = 99000
$ws.Cells.Item(90, 14).Value = -111480
$ws.Cells.Item(116, 8).Value = 728.5833
$ws.Cells.Item(116, 9).Value = 660.3333
$ws.Cells.Item(116, 11).Value = 660.3333
$ws.Cells.Item(116, 13).Value = 1633.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 728.5833
$ws.Cells.Item(3, 9).Value = 660.3333
$ws.Cells.Item(3, 11).Value = 660.3333
$ws.Cells.Item(3, 13).Value = -546.3333
$ws.Cells.Item(105, 8).Value = 2483.5857
$ws.Cells.Item(105, 9).Value = 2320.1924
$ws.Cells.Item(105, 11).Value = 2320.1924
$ws.Cells.Item(105, 13).Value = -573.1923999999999
$ws.Cells.Item(134, 8).Value = 23089.102
$ws.Cells.Item(134, 9).Value = 30707.824
$ws.Cells.Item(134, 10).Value = 5820
$ws.Cells.Item(134, 11).Value = 92123.47200000001
$ws.Cells.Item(134, 12).Value = 17460
$ws.Cells.Item(134, 13).Value = -89588.47200000001
$ws.Cells.Item(134, 14).Value = -22530

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 333933.34
$ws.Cells.Item(16, 9).Value = 333933.34
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 333933.34
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -333646.34
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 333933.34
$ws.Cells.Item(113, 9).Value = 333933.34
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 333933.34
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -331763.34
$ws.Cells.Item(113, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 503.05
$ws.Cells.Item(5, 10).Value = 870.9091
$ws.Cells.Item(5, 12).Value = 2612.7273
$ws.Cells.Item(5, 14).Value = -2836.7273
$ws.Cells.Item(122, 8).Value = 520.53125
$ws.Cells.Item(122, 9).Value = 384.31818
$ws.Cells.Item(122, 11).Value = 3458.86362
$ws.Cells.Item(122, 13).Value = -1008.86362
$ws.Cells.Item(132, 8).Value = 3814.3235
$ws.Cells.Item(132, 9).Value = 1361.9166
$ws.Cells.Item(132, 11).Value = 12257.2494
$ws.Cells.Item(132, 13).Value = -9727.249400000001
$ws.Cells.Item(135, 8).Value = 503.05
$ws.Cells.Item(135, 10).Value = 870.9091
$ws.Cells.Item(135, 12).Value = 7838.1819
$ws.Cells.Item(135, 14).Value = -12908.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9631.727999999999
$ws.Cells.Item(70, 9).Value = 3992.8572
$ws.Cells.Item(70, 11).Value = 3992.8572
$ws.Cells.Item(70, 13).Value = -3722.8572
$ws.Cells.Item(73, 8).Value = 9631.727999999999
$ws.Cells.Item(73, 9).Value = 3992.8572
$ws.Cells.Item(73, 11).Value = 3992.8572
$ws.Cells.Item(73, 13).Value = -3056.8572
$ws.Cells.Item(80, 8).Value = 4703.7744
$ws.Cells.Item(80, 9).Value = 5258.5415
$ws.Cells.Item(80, 10).Value = 2801.7144
$ws.Cells.Item(80, 11).Value = 5258.5415
$ws.Cells.Item(80, 12).Value = 2801.7144
$ws.Cells.Item(80, 13).Value = -4260.5415
$ws.Cells.Item(80, 14).Value = -4797.7144
$ws.Cells.Item(83, 8).Value = 4703.7744
$ws.Cells.Item(83, 9).Value = 5258.5415
$ws.Cells.Item(83, 10).Value = 2801.7144
$ws.Cells.Item(83, 11).Value = 26292.7075
$ws.Cells.Item(83, 12).Value = 14008.572
$ws.Cells.Item(83, 13).Value = -21300.7075
$ws.Cells.Item(83, 14).Value = -23992.572
$ws.Cells.Item(126, 8).Value = 2499.9443
$ws.Cells.Item(126, 9).Value = 2691.1538
$ws.Cells.Item(126, 10).Value = 2002.8
$ws.Cells.Item(126, 11).Value = 8073.4614
$ws.Cells.Item(126, 12).Value = 6008.4
$ws.Cells.Item(126, 13).Value = -5603.4614
$ws.Cells.Item(126, 14).Value = -10948.4
$ws.Cells.Item(134, 8).Value = 26285.715
$ws.Cells.Item(134, 10).Value = 26285.715
$ws.Cells.Item(134, 12).Value = 78857.145
$ws.Cells.Item(134, 14).Value = -83927.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2275.5
$ws.Cells.Item(7, 9).Value = 2286.5557
$ws.Cells.Item(7, 10).Value = 2261.2856
$ws.Cells.Item(7, 11).Value = 2286.5557
$ws.Cells.Item(7, 12).Value = 2261.2856
$ws.Cells.Item(7, 13).Value = -2174.5557
$ws.Cells.Item(7, 14).Value = -2485.2856
$ws.Cells.Item(82, 8).Value = 1703.5483
$ws.Cells.Item(82, 9).Value = 1355.05
$ws.Cells.Item(82, 10).Value = 2337.182
$ws.Cells.Item(82, 11).Value = 1355.05
$ws.Cells.Item(82, 12).Value = 2337.182
$ws.Cells.Item(82, 13).Value = -994.05
$ws.Cells.Item(82, 14).Value = -3059.182
$ws.Cells.Item(85, 8).Value = 1703.5483
$ws.Cells.Item(85, 9).Value = 1355.05
$ws.Cells.Item(85, 10).Value = 2337.182
$ws.Cells.Item(85, 11).Value = 1355.05
$ws.Cells.Item(85, 12).Value = 2337.182
$ws.Cells.Item(85, 13).Value = -107.05
$ws.Cells.Item(85, 14).Value = -4833.182
$ws.Cells.Item(87, 8).Value = 33729.668
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 33729.668
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 33729.668
$ws.Cells.Item(87, 13).ClearContents()
$ws.Cells.Item(87, 14).Value = -35975.668
$ws.Cells.Item(90, 8).Value = 33729.668
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 33729.668
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 101189.004
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(90, 14).Value = -112421.004
$ws.Cells.Item(126, 8).Value = 2275.5
$ws.Cells.Item(126, 9).Value = 2286.5557
$ws.Cells.Item(126, 10).Value = 2261.2856
$ws.Cells.Item(126, 11).Value = 6859.6671
$ws.Cells.Item(126, 12).Value = 6783.8568
$ws.Cells.Item(126, 13).Value = -4389.6671
$ws.Cells.Item(126, 14).Value = -11723.8568
$ws.Cells.Item(136, 8).Value = 4290.2
$ws.Cells.Item(136, 9).Value = 2407.6296
$ws.Cells.Item(136, 10).Value = 6500.174
$ws.Cells.Item(136, 11).Value = 7222.888800000001
$ws.Cells.Item(136, 12).Value = 19500.522
$ws.Cells.Item(136, 13).Value = -4672.888800000001
$ws.Cells.Item(136, 14).Value = -24600.522

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 11633.333
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 16950
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 12).Value = 16950
$ws.Cells.Item(4, 13).Value = -887
$ws.Cells.Item(4, 14).Value = -17176
$ws.Cells.Item(62, 8).Value = 3358.8235
$ws.Cells.Item(62, 9).Value = 2302
$ws.Cells.Item(62, 10).Value = 4193.1577
$ws.Cells.Item(62, 11).Value = 2302
$ws.Cells.Item(62, 12).Value = 4193.1577
$ws.Cells.Item(62, 13).Value = -1678
$ws.Cells.Item(62, 14).Value = -5441.1577
$ws.Cells.Item(65, 8).Value = 3358.8235
$ws.Cells.Item(65, 9).Value = 2302
$ws.Cells.Item(65, 10).Value = 4193.1577
$ws.Cells.Item(65, 11).Value = 11510
$ws.Cells.Item(65, 12).Value = 20965.7885
$ws.Cells.Item(65, 13).Value = -8390
$ws.Cells.Item(65, 14).Value = -27205.7885
$ws.Cells.Item(81, 8).Value = 2921
$ws.Cells.Item(81, 9).Value = 2181.3
$ws.Cells.Item(81, 11).Value = 4362.6
$ws.Cells.Item(81, 13).Value = -3301.6
$ws.Cells.Item(84, 8).Value = 2921
$ws.Cells.Item(84, 9).Value = 2181.3
$ws.Cells.Item(84, 11).Value = 21813
$ws.Cells.Item(84, 13).Value = -16509
$ws.Cells.Item(122, 8).Value = 35734.414
$ws.Cells.Item(122, 9).Value = 39620.31
$ws.Cells.Item(122, 10).Value = 2056.6667
$ws.Cells.Item(122, 11).Value = 118860.93
$ws.Cells.Item(122, 12).Value = 6170.000100000001
$ws.Cells.Item(122, 13).Value = -116410.93
$ws.Cells.Item(122, 14).Value = -11070.0001
$ws.Cells.Item(136, 8).Value = 2404.238
$ws.Cells.Item(136, 9).Value = 1961.4736
$ws.Cells.Item(136, 10).Value = 2770
$ws.Cells.Item(136, 12).Value = 8310
$ws.Cells.Item(136, 13).Value = -3334.4208
$ws.Cells.Item(136, 14).Value = -13410
